$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-17 from 45184 (2023-09-15)
# to 45185 (2023-09-16), keeping the existing date formatting.
$ws.Range("C2:C17").Value = 45185
